# Reverse the order of comma-separated "Recorded By" entries in column G.
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Cells with only a single (non comma-separated) entry are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G holds the "Recorded By" values
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value()

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        $n = $parts.Count
        $result = ""
        for ($i = $n - 1; $i -ge 0; $i--) {
            $piece = $parts[$i].Trim()
            if ($result -eq "") {
                $result = $piece
            } else {
                $result = $result + ", " + $piece
            }
        }
        $cell.Value = $result
    }
}
